# Generate Report for Handback
#
# The "32621205-ad4a-4fb9-ab0b-ccf4eba2277d" file has finished its
# handback (it used to be "Ready for handoff"); its row now reports
# "Handed back: in sync with en-US" just like the "a1a62ee8-..." file,
# the two rows swap display order (32621205 now listed first), and the
# handback timestamps for both locales are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-22 07:09:38"

$ws1.Range("A3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-22 07:09:38"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6141f81a70034703285afa96d6b965b17ebd869/e2e/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a5e6915cbd0b22f9f71240ceed01ff4d2acd7a0b/e2e/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-22 07:09:29"
$ws2.Range("F2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md"
$ws2.Range("G2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-22 07:10:16"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-22 07:09:29"
$ws2.Range("F3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md"
$ws2.Range("G3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-22 07:10:16"
$ws2.Range("J3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6141f81a70034703285afa96d6b965b17ebd869/e2e/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/670273ed977333b1f2b31fd7c35cbf37dbbd8d60/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/mt/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.zh-cn.xlf", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/e99f04c9eca89f780d05e594ed238e4022fc3365/e2e/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/773e2841639e85edbbc866af6adeb9f56760243b/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.zh-cn.xlf", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a5e6915cbd0b22f9f71240ceed01ff4d2acd7a0b/e2e/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/670273ed977333b1f2b31fd7c35cbf37dbbd8d60/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/mt/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.zh-cn.xlf", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/e99f04c9eca89f780d05e594ed238e4022fc3365/e2e/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/773e2841639e85edbbc866af6adeb9f56760243b/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.zh-cn.xlf", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-22 07:09:38"
$ws3.Range("F2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md"
$ws3.Range("G2").Value = "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-22 07:10:35"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-22 07:09:38"
$ws3.Range("F3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md"
$ws3.Range("G3").Value = "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-22 07:10:35"
$ws3.Range("J3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6141f81a70034703285afa96d6b965b17ebd869/e2e/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1c48f226591b93ed8577b8ea5654b2a1264abbb/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/mt/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.de-de.xlf", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/f37f496982ff66c2065581d24f23385680fd300c/e2e/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8494097df150b6d25b5675095929118fd9f57fb3/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.de-de.xlf", "", "", "32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a5e6915cbd0b22f9f71240ceed01ff4d2acd7a0b/e2e/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1c48f226591b93ed8577b8ea5654b2a1264abbb/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/mt/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.de-de.xlf", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/f37f496982ff66c2065581d24f23385680fd300c/e2e/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.md", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8494097df150b6d25b5675095929118fd9f57fb3/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/32621205-ad4a-4fb9-ab0b-ccf4eba2277d.d4353c2348a121f89e93a17cc55e405fb5837e0b.de-de.xlf", "", "", "a1a62ee8-c0ae-4d85-90a1-c93fd22bc574.200ce595132812212e5eda58ac0a0447321cb687.de-de.xlf")
